$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $addr, $val) {
    $cell = $ws.Range($addr)
    # Force a text number format before writing so Excel does not
    # auto-convert numeric-looking strings (e.g. "1.00" -> 1, dropping
    # the trailing zeros / decimal points used as thousands separators).
    $cell.NumberFormat = "@"
    $cell.Value = $val
    # Restore the default "Normal" style so the cell's style index (s attribute)
    # is unchanged from before the edit - only the cell content changes.
    $cell.Style = "Normal"
}

Set-TextValue $ws 'D2' '60.514.50'
Set-TextValue $ws 'E2' '  -2.91%  '
Set-TextValue $ws 'D3' '2.905.36'
Set-TextValue $ws 'E3' '  -3.38%  '
Set-TextValue $ws 'E4' '  +0.04%  '
Set-TextValue $ws 'D5' '526.69'
Set-TextValue $ws 'E5' '  -5.08%  '
Set-TextValue $ws 'D6' '141.43'
Set-TextValue $ws 'E6' '  -7.33%  '
Set-TextValue $ws 'D7' '1.00'
Set-TextValue $ws 'E7' '  -0.03%  '
Set-TextValue $ws 'D8' '0.550'
Set-TextValue $ws 'E8' '  -3.52%  '
Set-TextValue $ws 'D9' '2.911.99'
Set-TextValue $ws 'E9' '  -3.28%  '
Set-TextValue $ws 'E10' '  -5.77%  '
Set-TextValue $ws 'D11' '5.90'
Set-TextValue $ws 'E11' '  -6.21%  '
Set-TextValue $ws 'D12' '0.355'
Set-TextValue $ws 'E12' '  -3.44%  '
Set-TextValue $ws 'D13' '3.417.82'
Set-TextValue $ws 'E13' '  -3.22%  '
Set-TextValue $ws 'E14' '  +1.12%  '
Set-TextValue $ws 'D15' '60.549.86'
Set-TextValue $ws 'D16' '22.64'
Set-TextValue $ws 'E16' '  -4.94%  '
Set-TextValue $ws 'D17' '2.910.20'
Set-TextValue $ws 'E17' '  -3.20%  '
Set-TextValue $ws 'E18' '  -6.36%  '
Set-TextValue $ws 'D19' '4.93'
Set-TextValue $ws 'E19' '  -3.72%  '
Set-TextValue $ws 'D20' '11.52'
Set-TextValue $ws 'E20' '  -4.18%  '
Set-TextValue $ws 'D21' '361.38'
Set-TextValue $ws 'E21' '  -8.37%  '
Set-TextValue $ws 'D22' '6.56'
Set-TextValue $ws 'E22' '  -2.20%  '
Set-TextValue $ws 'D23' '0.999'
Set-TextValue $ws 'E23' '  -0.21%  '
Set-TextValue $ws 'D24' '63.33'
Set-TextValue $ws 'E24' '  -2.82%  '
Set-TextValue $ws 'D25' '3.018.79'
Set-TextValue $ws 'E25' '  -3.89%  '
Set-TextValue $ws 'D26' '0.447'
Set-TextValue $ws 'E26' '  -4.80%  '
Set-TextValue $ws 'D27' '0.180'
Set-TextValue $ws 'E27' '  -3.69%  '
Set-TextValue $ws 'D28' '0.999'
Set-TextValue $ws 'E28' '  -0.31%  '
Set-TextValue $ws 'D29' '7.79'
Set-TextValue $ws 'E29' '  -8.80%  '
Set-TextValue $ws 'D30' '0.0₃0848'
Set-TextValue $ws 'E30' '  -12.91%  '
Set-TextValue $ws 'D31' '0.999'
Set-TextValue $ws 'E31' '  -0.02%  '
Set-TextValue $ws 'E32' '  -4.65%  '
Set-TextValue $ws 'D33' '19.46'
Set-TextValue $ws 'E33' '  -5.38%  '
Set-TextValue $ws 'D34' '150.90'
Set-TextValue $ws 'E34' '  -5.98%  '
Set-TextValue $ws 'D35' '4.30'
Set-TextValue $ws 'E35' '  -8.74%  '
Set-TextValue $ws 'D36' '5.52'
Set-TextValue $ws 'E36' '  -8.78%  '
Set-TextValue $ws 'E37' '  -9.45%  '
Set-TextValue $ws 'E38' '  -8.55%  '
Set-TextValue $ws 'D39' '37.94'
Set-TextValue $ws 'E39' '  +1.09%  '
Set-TextValue $ws 'E40' '  -6.08%  '
Set-TextValue $ws 'D41' '2.330.28'
Set-TextValue $ws 'E41' '  -5.27%  '
Set-TextValue $ws 'D42' '0.645'
Set-TextValue $ws 'E42' '  -2.73%  '
Set-TextValue $ws 'E43' '  -7.48%  '
Set-TextValue $ws 'D44' '20.70'
Set-TextValue $ws 'E44' '  -8.44%  '
Set-TextValue $ws 'D45' '0.0570'
Set-TextValue $ws 'E45' '  -4.78%  '
Set-TextValue $ws 'E46' '  -0.05%  '
Set-TextValue $ws 'E47' '  +0.23%  '
Set-TextValue $ws 'D48' '0.0233'
Set-TextValue $ws 'E48' '  -6.38%  '
Set-TextValue $ws 'D49' '10.31'
Set-TextValue $ws 'E49' '  -1.75%  '
Set-TextValue $ws 'D50' '0.0926'
Set-TextValue $ws 'E50' '  -3.12%  '
Set-TextValue $ws 'D51' '250.85'
Set-TextValue $ws 'E51' '  -5.29%  '
